$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "30÷2=" "20÷6="
Replace-Text "91÷8=" "71÷7="
Replace-Text "47÷8=" "52÷4="
Replace-Text "91÷6=" "22÷3="
Replace-Text "45÷8=" "78÷4="
Replace-Text "87÷4=" "57÷6="
Replace-Text "22÷9=" "31÷4="
Replace-Text "82÷4=" "45÷3="
Replace-Text "68÷2=" "55÷2="
Replace-Text "12÷5=" "48÷8="
Replace-Text "88÷5=" "74÷6="
Replace-Text "24÷3=" "98÷9="
Replace-Text "63÷9=" "51÷4="
Replace-Text "51÷8=" "81÷7="
Replace-Text "50÷4=" "32÷4="
Replace-Text "66÷9=" "46÷9="
Replace-Text "90÷6=" "86÷4="
Replace-Text "50÷8=" "38÷6="
Replace-Text "93÷2=" "68÷8="
Replace-Text "96÷9=" "84÷5="
Replace-Text "76÷7=" "73÷3="
Replace-Text "91÷2=" "25÷6="
Replace-Text "99÷6=" "63÷3="
Replace-Text "54÷3=" "46÷2="
Replace-Text "48÷4=" "24÷3="
